# Add a new worksheet "8_" (a coffee-cooling / deltaT quiz page) at the end
# of the workbook, built as a copy of the last existing sheet "7_" so that
# column widths, row-label styles (s="1" wrap / s="2" centered) and overall
# layout are inherited automatically, then overwrite its question content.

$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("7_")

# Copy "7_" and place the duplicate immediately after it. Copy() duplicates
# sheet formatting/col widths/row heights and also moves the "active tab"
# marker onto the new sheet (matching tabSelected moving off of "7_").
$src.Copy([System.Reflection.Missing]::Value, $src)
$new = $wb.Worksheets.Item($wb.Worksheets.Count)
$new.Name = "8_"

# --- Row 1: question prompt / column headers ---
$new.Range("A1").Value = "Look at the discretized equation for deltaT.  Which factors would cause the coffee to cool more quickly over the length of the simulation?"
$new.Range("B1").Value = "Correct order of definitions"
$new.Range("C1").Value = "Definitions"

# --- Row 2: "A larger time step dt" -> B (incorrect) ---
$new.Range("A2").Value = "A larger time step dt"
$new.Range("B2").Value = "B"
$new.Range("C2").Value = "This will cause the coffee to cool more quickly"
$new.Range("D2").Value = "This is a bit of a trick.  While a large dt will cause a bigger change in a single time step, there will be fewer time steps in the whole simulation.  dt should not significantly affect the cooling rate (until it gets much too big)."

# --- Row 3: "A higher T_env" -> B (incorrect) ---
$new.Range("A3").Value = "A higher T_env"
$new.Range("B3").Value = "B"
$new.Range("C3").Value = "This will not cause the coffee to cool more quickly"
$new.Range("D3").Value = "This will lower the temperature difference between the air and the coffee, and the coffee will cool slower"

# --- Row 4: "A higher initial T_init" -> A (correct) ---
$new.Range("A4").Value = "A higher initial T_init"
$new.Range("B4").Value = "A"
$new.Range("C4").Clear()
$new.Range("D4").Clear()

# --- Row 5: "A larger r" -> A (correct) ---
$new.Range("A5").Value = "A larger r"
$new.Range("B5").Value = "A"
$new.Range("C5").Clear()
$new.Range("D5").Value = "Factors that might lead to a larger r include more convective cooling or poorer insulation."

# --- Rows 6-7: leftover answer-entry cells from the copied sheet are blanked ---
$new.Range("A6").Clear()
$new.Range("B6").ClearContents()
$new.Range("D6").Clear()

# Row heights: 90 / 120 / 60 / (default) / 45 / (default) / (default)
$new.Rows.Item(1).RowHeight = 90
$new.Rows.Item(2).RowHeight = 120
$new.Rows.Item(3).RowHeight = 60
$new.Rows.Item(4).AutoFit()
$new.Rows.Item(5).RowHeight = 45
$new.Rows.Item(6).AutoFit()

# Leave the selection on the new sheet's blank answer area, same as authored.
$new.Range("F5:F7").Select() | Out-Null
